$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The "Metadata" sheet had a duplicated "Contact" / "No display for ContactDetail"
# row (rows 10 and 11 were identical). Remove the second (duplicate) one - this
# shifts everything below it up by one row.
$ws1.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$ws1.Cells.Item(3, 2).Value = "6.0.0"

# Date: refreshed publication date
$ws1.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Publisher value was empty, now populated
$ws1.Cells.Item(9, 2).Value = "Alvearie Team"

# The remaining "Contact" row becomes a "Jurisdiction" row
$ws1.Cells.Item(10, 1).Value = "Jurisdiction"
$ws1.Cells.Item(10, 2).Value = "United States of America"

# "Case Sensitive" value was empty, now "true". Assigning a literal "true"
# through .Value auto-coerces to a Boolean, which would change the cell's
# stored type (and force a new number-format style) versus the plain text
# shared-string the source file uses. Route it through a text formula and
# then flatten the formula to a literal value so it lands as text "true".
$caseSensitiveCell = $ws1.Cells.Item(14, 2)
$caseSensitiveCell.Formula = '="true"'
$caseSensitiveCell.Copy()
$caseSensitiveCell.PasteSpecial(-4163)
